$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Row 37: 2C_Metal-production activity data reverted to the standard
# GDP / B2005USD pattern used by the rest of the sheet, removing the
# older "Metal_S_Content" smelting activity + its description.
$ws.Range("B37").Value = "GDP"
$ws.Range("C37").Value = "B2005USD"
$ws.Range("D37").Value = ""

# Row 42: 2H_Pulp-and-paper-food-beverage-wood activity data reverted to
# the standard GDP / B2005USD pattern, removing the older
# "Pulp_Paper_Production" / "tons" activity data.
$ws.Range("B42").Value = "GDP"
$ws.Range("C42").Value = "B2005USD"

# Reflect the author's final on-screen selection/scroll position in the
# saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("B30").Select()
